$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header row swap: Foliar Al/Ca/P rotate between columns C, E, K ---
$ws.Range("C1").Value2 = "Foliar Al"
$ws.Range("E1").Value2 = "Foliar Ca"
$ws.Range("K1").Value2 = "Foliar P"

# --- 2. Rename the stat-column headers (row 2): Chisq -> F value, Pr(>Chisq) -> Pr(>F) ---
$ws.Cells.Replace("Chisq", "F value", 1, 1, $false, $false, $false, $false) | Out-Null
$ws.Cells.Replace("Pr(>Chisq)", "Pr(>F)", 1, 1, $false, $false, $false, $false) | Out-Null

# --- 3. New linear-model results for rows 3-5 (columns C:N) ---
$ws.Range("C3").Value2 = 0.393432052299655
$ws.Range("D3").Value2 = 0.5344591491596
$ws.Range("E3").Value2 = 15.3955436822811
$ws.Range("F3").Value2 = 0.00037643470054839
$ws.Range("G3").Value2 = 5.0421541853429
$ws.Range("H3").Value2 = 0.0309714682930822
$ws.Range("I3").Value2 = 2.20080650946547
$ws.Range("J3").Value2 = 0.146641524218118
$ws.Range("K3").Value2 = 0.523060704836498
$ws.Range("L3").Value2 = 0.474212121152692
$ws.Range("M3").Value2 = 8.64480114827674
$ws.Range("N3").Value2 = 0.00569944021193973

$ws.Range("C4").Value2 = 0.0206622884549123
$ws.Range("D4").Value2 = 0.886504397925486
$ws.Range("E4").Value2 = 0.84346711131403
$ws.Range("F4").Value2 = 0.36451967354065
$ws.Range("G4").Value2 = 4.07130221938722
$ws.Range("H4").Value2 = 0.0511187651515683
$ws.Range("I4").Value2 = 0.507436514949917
$ws.Range("J4").Value2 = 0.480842281829567
$ws.Range("K4").Value2 = 8.30877468245125
$ws.Range("L4").Value2 = 0.00661518524626132
$ws.Range("M4").Value2 = 0.0499171172918882
$ws.Range("N4").Value2 = 0.824470552190544

$ws.Range("C5").Value2 = 0.187329294967568
$ws.Range("D5").Value2 = 0.667730820769702
$ws.Range("E5").Value2 = 0.0877810934236624
$ws.Range("F5").Value2 = 0.768719743413029
$ws.Range("G5").Value2 = 4.86313289395857
$ws.Range("H5").Value2 = 0.033908823002517297
$ws.Range("I5").Value2 = 0.377388099590421
$ws.Range("J5").Value2 = 0.542867619451516
$ws.Range("K5").Value2 = 0.407348862067539
$ws.Range("L5").Value2 = 0.527357672965897
$ws.Range("M5").Value2 = 1.45767660153883
$ws.Range("N5").Value2 = 0.23517307257084

# --- 4. Drop the leftover highlight style on B3:B5 and C2:N2 (back to the sheet default) ---
$ws.Range("C2:N2").Style = "Normal"
$ws.Range("B3").Style = "Normal"
$ws.Range("B4").Style = "Normal"
$ws.Range("B5").Style = "Normal"

# --- 5. New trailing helper cells that ride along with the refreshed ANOVA table/paste ---
foreach ($addr in @("Q3","Q4","Q5")) {
  $ws.Range($addr).Borders.LineStyle = 1
  $ws.Range($addr).Borders.LineStyle = 0
}
$ws.Range("C6:S6").Borders.LineStyle = 1
$ws.Range("C6:S6").Borders.LineStyle = 0
$ws.Range("B6").Value2 = 36
$ws.Range("B6").Style = "Normal"

# --- 6. Column widths re-fit to the new (narrower) numbers ---
$ws.Columns.Item(2).ColumnWidth = 2.285714285714286
$ws.Columns.Item(3).ColumnWidth = 7.142857142857143
$ws.Columns.Item(4).ColumnWidth = 8.714285714285714
$ws.Columns.Item(5).ColumnWidth = 7.428571428571429
$ws.Columns.Item(6).ColumnWidth = 8.714285714285714
$ws.Columns.Item(7).ColumnWidth = 6.428571428571429
$ws.Columns.Item(8).ColumnWidth = 8.714285714285714
$ws.Columns.Item(9).ColumnWidth = 8.142857142857142
$ws.Columns.Item(10).ColumnWidth = 8.714285714285714
$ws.Columns.Item(11).ColumnWidth = 6.428571428571429
$ws.Columns.Item(12).ColumnWidth = 8.714285714285714
$ws.Columns.Item(13).ColumnWidth = 7.428571428571429
$ws.Columns.Item(14).ColumnWidth = 8.714285714285714

# --- 7. Selection moved to cover the (now one row taller) results block ---
$ws.Range("B3:N6").Select() | Out-Null
